# "wrapping up test file audit"
#
# The optimization_parameters sheet had a stray leftover row labeled
# "Sheet" (with junk values 3 / 4) sitting between the "Strain" row and
# the "simulation_timepoints" row. Remove it — this also shifts the
# "simulation_timepoints" row up by one and, because "Sheet" was only
# referenced by that one cell, drops it out of the shared-string table.
$wb = $excel.ActiveWorkbook

$paramsSheet = $wb.Worksheets.Item("optimization_parameters")
$paramsSheet.Rows("16:16").Delete()

# The audit pass finished on the optimization_diagnostics tab, so that's
# the sheet left active/selected when the workbook was saved.
$diagSheet = $wb.Worksheets.Item("optimization_diagnostics")
$diagSheet.Activate()
